$d = $word.ActiveDocument

$d.Content.Find.Execute("Item 7: " + [char]0x2026, $true, $false, $false, $false, $false, $true, 1, $false, "Item 7: Tutorial", 2)

$d.Content.Find.Execute("Tipo de melhoria: " + [char]0x2026 + "  ", $true, $false, $false, $false, $false, $true, 1, $false, "Tipo de melhoria: Tutorial para a primeira jogatina do usuário.  ", 2)

$d.Content.Find.Execute("Estado anterior:... ", $true, $false, $false, $false, $false, $true, 1, $false, "Estado anterior: Funcionalidades sem indicação.", 2)

$d.Content.Find.Execute("Descrição da melhoria: " + [char]0x2026, $true, $false, $false, $false, $false, $true, 1, $false, "Descrição da melhoria: Alguns pop ups indicando como funciona a mecânica geral do jogo.", 2)

$d.Content.Find.Execute("Foi implementada?...", $true, $false, $false, $false, $false, $true, 1, $false, "Foi implementada?: Ainda não.", 2)
